$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "302.57"
Set-TextValue $ws.Range("E2") "0.79%"
Set-TextValue $ws.Range("G2") "20"
Set-TextValue $ws.Range("D3") "32.10"
Set-TextValue $ws.Range("E3") "0.91%"
Set-TextValue $ws.Range("G3") "20"
Set-TextValue $ws.Range("D4") "4.987"
Set-TextValue $ws.Range("E4") "-2.40%"
Set-TextValue $ws.Range("G4") "20"
Set-TextValue $ws.Range("D5") "0.07903"
Set-TextValue $ws.Range("E5") "-3.81%"
Set-TextValue $ws.Range("G5") "20"
Set-TextValue $ws.Range("D6") "2.161"
Set-TextValue $ws.Range("E6") "-17.64%"
Set-TextValue $ws.Range("G6") "20"
Set-TextValue $ws.Range("D7") "7.863"
Set-TextValue $ws.Range("E7") "0.17%"
Set-TextValue $ws.Range("G7") "20"
Set-TextValue $ws.Range("D8") "3.797"
Set-TextValue $ws.Range("G8") "20"
Set-TextValue $ws.Range("E9") "0.01%"
Set-TextValue $ws.Range("G9") "20"
Set-TextValue $ws.Range("D10") "0.1751"
Set-TextValue $ws.Range("E10") "-0.24%"
Set-TextValue $ws.Range("G10") "20"
Set-TextValue $ws.Range("D11") "0.08035"
Set-TextValue $ws.Range("E11") "6.95%"
Set-TextValue $ws.Range("G11") "20"
Set-TextValue $ws.Range("D12") "0.08780"
Set-TextValue $ws.Range("E12") "-3.33%"
Set-TextValue $ws.Range("G12") "20"
Set-TextValue $ws.Range("D13") "0.03132"
Set-TextValue $ws.Range("E13") "3.54%"
Set-TextValue $ws.Range("G13") "20"
Set-TextValue $ws.Range("E14") "0.19%"
Set-TextValue $ws.Range("G14") "20"
Set-TextValue $ws.Range("D15") "0.001514"
Set-TextValue $ws.Range("E15") "0.38%"
Set-TextValue $ws.Range("G15") "20"
Set-TextValue $ws.Range("D16") "0.005970"
Set-TextValue $ws.Range("E16") "-1.32%"
Set-TextValue $ws.Range("G16") "20"
Set-TextValue $ws.Range("E17") "-4.03%"
Set-TextValue $ws.Range("G17") "20"
Set-TextValue $ws.Range("D18") "2.280"
Set-TextValue $ws.Range("E18") "-0.24%"
Set-TextValue $ws.Range("G18") "20"
Set-TextValue $ws.Range("G19") "20"
Set-TextValue $ws.Range("D20") "0.1290"
Set-TextValue $ws.Range("E20") "-4.14%"
Set-TextValue $ws.Range("G20") "20"
Set-TextValue $ws.Range("D21") "4.146"
Set-TextValue $ws.Range("E21") "5.87%"
Set-TextValue $ws.Range("G21") "20"
Set-TextValue $ws.Range("D22") "0.1791"
Set-TextValue $ws.Range("E22") "6.79%"
Set-TextValue $ws.Range("G22") "20"
Set-TextValue $ws.Range("D23") "0.04606"
Set-TextValue $ws.Range("E23") "-0.18%"
Set-TextValue $ws.Range("G23") "20"
Set-TextValue $ws.Range("D24") "0.001237"
Set-TextValue $ws.Range("E24") "-0.66%"
Set-TextValue $ws.Range("G24") "20"
Set-TextValue $ws.Range("D25") "0.004503"
Set-TextValue $ws.Range("E25") "-1.18%"
Set-TextValue $ws.Range("G25") "20"
Set-TextValue $ws.Range("E26") "4.23%"
Set-TextValue $ws.Range("G26") "20"
Set-TextValue $ws.Range("G27") "20"
Set-TextValue $ws.Range("G28") "20"
Set-TextValue $ws.Range("G29") "20"
Set-TextValue $ws.Range("G30") "20"
Set-TextValue $ws.Range("G31") "20"
Set-TextValue $ws.Range("G32") "20"
Set-TextValue $ws.Range("G33") "20"
Set-TextValue $ws.Range("G34") "20"
Set-TextValue $ws.Range("G35") "20"
Set-TextValue $ws.Range("G36") "20"
Set-TextValue $ws.Range("G37") "20"
Set-TextValue $ws.Range("G38") "20"
Set-TextValue $ws.Range("D39") "0.01730"
Set-TextValue $ws.Range("E39") "-2.72%"
Set-TextValue $ws.Range("G39") "20"
Set-TextValue $ws.Range("D40") "0.04830"
Set-TextValue $ws.Range("E40") "4.74%"
Set-TextValue $ws.Range("G40") "20"
Set-TextValue $ws.Range("D41") "0.007341"
Set-TextValue $ws.Range("E41") "6.70%"
Set-TextValue $ws.Range("G41") "20"
Set-TextValue $ws.Range("D42") "0.1367"
Set-TextValue $ws.Range("E42") "-1.05%"
Set-TextValue $ws.Range("G42") "20"
Set-TextValue $ws.Range("D43") "0.002338"
Set-TextValue $ws.Range("E43") "5.94%"
Set-TextValue $ws.Range("G43") "20"
Set-TextValue $ws.Range("D44") "0.01102"
Set-TextValue $ws.Range("E44") "12.94%"
Set-TextValue $ws.Range("G44") "20"
Set-TextValue $ws.Range("D45") "0.00006008"
Set-TextValue $ws.Range("E45") "-2.48%"
Set-TextValue $ws.Range("G45") "20"
Set-TextValue $ws.Range("E46") "0.21%"
Set-TextValue $ws.Range("G46") "20"
Set-TextValue $ws.Range("D47") "0.003388"
Set-TextValue $ws.Range("E47") "-59.56%"
Set-TextValue $ws.Range("G47") "20"
Set-TextValue $ws.Range("D48") "0.8234"
Set-TextValue $ws.Range("E48") "4.23%"
Set-TextValue $ws.Range("G48") "20"
Set-TextValue $ws.Range("E49") "0.21%"
Set-TextValue $ws.Range("G49") "20"
Set-TextValue $ws.Range("E50") "0.21%"
Set-TextValue $ws.Range("G50") "20"
Set-TextValue $ws.Range("G51") "20"
